$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.272.20'
$ws.Range("E2").Value = '  +3.52%  '
$ws.Range("D3").Value = '3.484.14'
$ws.Range("E3").Value = '  +2.58%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.47%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +12.08%  '
$ws.Range("D9").Value = '3.485.84'
$ws.Range("E9").Value = '  +2.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("E11").Value = '  +3.38%  '
$ws.Range("E12").Value = '  +4.07%  '
$ws.Range("D13").Value = '4.088.36'
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("E15").Value = '  +3.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.64%  '
$ws.Range("D17").Value = '65.328.46'
$ws.Range("E17").Value = '  +3.42%  '
$ws.Range("D18").Value = '3.469.62'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '387.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("E23").Value = '  +4.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  +5.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.76%  '
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +11.30%  '
$ws.Range("E31").Value = '  +4.44%  '
$ws.Range("E32").Value = '  +3.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.11%  '
$ws.Range("E37").Value = '  +5.88%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '2.999.74'
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0778'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0325'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.77%  '
$ws.Range("E42").Value = '  +6.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.83%  '
$ws.Range("E44").Value = '  +3.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.779'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '322.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.61%  '
$ws.Range("E49").Value = '  +6.17%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.03%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.109'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.68%  '
